# Workbook was edited: stimuli pictures were removed from the Stimuli folder,
# and "Tabelle1" was turned into a readme-style list of the remaining image files.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

$values = @(
    "ImageFile",
    "Stimuli/223.jpg",
    "Stimuli/238.jpg",
    "Stimuli/245.jpg",
    "Stimuli/2981.jpg",
    "Stimuli/3016.jpg",
    "Stimuli/3101.jpg",
    "Stimuli/3181.jpg",
    "Stimuli/3215.jpg",
    "Stimuli/3220.jpg",
    "Stimuli/3225.jpg",
    "Stimuli/6020.jpg",
    "Stimuli/6571.jpg",
    "Stimuli/6831.jpg",
    "Stimuli/8231.jpg",
    "Stimuli/9373.jpg",
    "Stimuli/9402.jpg",
    "Stimuli/9400.jpg",
    "Stimuli/9403.jpg",
    "Stimuli/9405.jpg",
    "Stimuli/9423.jpg"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws1.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws1.Range("A2:A21").Select() | Out-Null
